# Apply the "SIQ-" ID prefix rename and the RATE-01 answer-text tweak,
# matching the commit "Add files via upload" re-upload of the
# SIQ Travel Advisor Web Application requirements sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the requirement IDs in column A (A2:A10) to the "SIQ-" scheme
# ---------------------------------------------------------------------
$idMap = @{
    "A2"  = "SIQ-Reg-01"
    "A3"  = "SIQ-Reg-02"
    "A4"  = "SIQ-Reg-03"
    "A5"  = "SIQ-Login-01"
    "A6"  = "SIQ-Admin-01"
    "A7"  = "SIQ-Rate-01"
    "A8"  = "SIQ-Booking-01"
    "A9"  = "SIQ-Gallery-01"
    "A10" = "SIQ-User-01"
}
foreach ($addr in $idMap.Keys) {
    $ws.Range($addr).Value = $idMap[$addr]
}

# ---------------------------------------------------------------------
# 2. Update the rating-system question (B7): "5-star scale" becomes
#    "5-star scale with text box", keeping the bold/plain rich-text runs.
# ---------------------------------------------------------------------
$cell = $ws.Range("B7")

# Touching the cell's whole-range font first mints the plain black Calibri
# style used as the cell's base format (matches the new cellXfs/font that
# the real edit introduced).
$cell.Font.Color = 0

$oldText = $cell.Characters().Text
$newText = $oldText.Replace("5-star scale", "5-star scale with text box")
$cell.Value = $newText
$text = $cell.Characters().Text

function Set-BoldRun($targetCell, $fullText, $substr) {
    $idx = $fullText.IndexOf($substr)
    if ($idx -ge 0) {
        $targetCell.Characters($idx + 1, $substr.Length).Font.Bold = $true
    }
}

function Set-ColorRun($targetCell, $fullText, $substr) {
    $idx = $fullText.IndexOf($substr)
    if ($idx -ge 0) {
        $targetCell.Characters($idx + 1, $substr.Length).Font.Color = 0
    }
}

Set-BoldRun $cell $text "rating system"
Set-BoldRun $cell $text "5-star scale with text box"
Set-BoldRun $cell $text "user history"

Set-ColorRun $cell $text "Should the "
Set-ColorRun $cell $text " use a "
Set-ColorRun $cell $text " and be displayed in "
Set-ColorRun $cell $text "?"

# ---------------------------------------------------------------------
# 3. Row-height touch-ups (auto-fit drift captured by the re-save)
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 28.9
$ws.Rows.Item(2).RowHeight = 45.75
$ws.Rows.Item(3).RowHeight = 45.75
$ws.Rows.Item(4).RowHeight = 60.75
$ws.Rows.Item(5).RowHeight = 30.75
$ws.Rows.Item(6).RowHeight = 45.75
$ws.Rows.Item(7).RowHeight = 45.75
$ws.Rows.Item(8).RowHeight = 45.75
$ws.Rows.Item(9).RowHeight = 45.75
$ws.Rows.Item(10).RowHeight = 45.75

# ---------------------------------------------------------------------
# 4. Restore the active selection left behind by the edit (C7)
# ---------------------------------------------------------------------
$ws.Range("C7").Select()
